# Apply "feat: add 2022-Q1 data" edit:
#  - The existing "总计" (Grand total) sheet (sheet index 3) becomes the new
#    "2022-Q1" per-fund holdings sheet.
#  - A brand-new "总计" (Grand total) sheet is appended at the end, containing
#    the original grand-total rows plus a new 2022-Q1 summary row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: turn the current "总计" sheet into the new "2022-Q1" sheet
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item(3)
$q1.Name = "2022-Q1"

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# B1 already carries the bordered/bold "header" style (s=2, inherited from the
# original sheet) -- copy that formatting onto the freshly-added header cells.
$q1.Range("B1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)

# Fund holdings detail rows:
#   rowNum, rowIndex, fundCode, fundName, fundSize, stockPosition, positionRatio, marketValue, positionRank
$fundRows = @(
    @(2, 0, "005119", "银华智荟内在价值灵活配置混合", "5.28", "94.60", "4.81", "0.2540", 9),
    @(3, 1, "009859", "银华乐享混合",                 "5.21", "94.63", "4.77", "0.2485", 9),
    @(4, 2, "005343", "长安裕盛灵活配置混合A",         "4.66", "94.21", "5.22", "0.2433", 7),
    @(5, 3, "005344", "长安裕盛灵活配置混合C",         "3.75", "94.21", "5.22", "0.1958", 7),
    @(6, 4, "004557", "北信瑞丰鼎丰灵活配置混合",       "0.39", "64.13", "5.00", "0.0195", 7)
)

foreach ($fr in $fundRows) {
    $r = $fr[0]

    $q1.Cells.Item($r, 1).Value = $fr[1]              # A: row index (number)

    $q1.Cells.Item($r, 2).Value = "'" + $fr[2]         # B: fund code (force text, keep leading zeros)
    $q1.Cells.Item($r, 2).Style = "Normal"

    $q1.Cells.Item($r, 3).Value = $fr[3]               # C: fund name (plain text)

    $q1.Cells.Item($r, 4).Value = "'" + $fr[4]         # D: fund size (force text)
    $q1.Cells.Item($r, 4).Style = "Normal"

    $q1.Cells.Item($r, 5).Value = "'" + $fr[5]         # E: stock position (force text)
    $q1.Cells.Item($r, 5).Style = "Normal"

    $q1.Cells.Item($r, 6).Value = "'" + $fr[6]         # F: position ratio (force text)
    $q1.Cells.Item($r, 6).Style = "Normal"

    $q1.Cells.Item($r, 7).Value = "'" + $fr[7]         # G: market value (force text, keep trailing zeros)
    $q1.Cells.Item($r, 7).Style = "Normal"

    $q1.Cells.Item($r, 8).Value = $fr[8]               # H: position rank (number)
}

# Column A (row index) also carries the s=2 style for rows 2..6; A2 already has
# it from the original sheet, so propagate it down to the newly-added rows.
$q1.Range("A2").Copy()
$q1.Range("A3:A6").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Step 2: append a brand-new "总计" sheet with the updated grand totals
# ---------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$total = $wb.Worksheets.Add($null, $lastSheet)
$total.Name = "总计"

# Match page margins used by the other data sheets (0.75in/0.75in/1in/1in/0.5in/0.5in)
$total.PageSetup.LeftMargin = 54
$total.PageSetup.RightMargin = 54
$total.PageSetup.TopMargin = 72
$total.PageSetup.BottomMargin = 72
$total.PageSetup.HeaderMargin = 36
$total.PageSetup.FooterMargin = 36

# Match outline settings used by the other data sheets
$total.Outline.SummaryRow = 1
$total.Outline.SummaryColumn = 1

# Header row
$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

# Copy the bordered/bold header formatting across from the 2022-Q1 sheet.
$q1.Range("B1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)

# Grand-total rows: rowNum, rowIndex, quarter, count, marketValue
$totalRows = @(
    @(2, 0, "2022-Q1", 5, 0.96),
    @(3, 1, "2021-Q4", 8, 1.11),
    @(4, 2, "2021-Q3", 13, 2.94)
)

foreach ($tr in $totalRows) {
    $r = $tr[0]
    $total.Cells.Item($r, 1).Value = $tr[1]   # A: row index (number)
    $total.Cells.Item($r, 2).Value = $tr[2]   # B: quarter label (plain text)
    $total.Cells.Item($r, 3).Value = $tr[3]   # C: fund count (number)
    $total.Cells.Item($r, 4).Value = $tr[4]   # D: market value (number)
}

# Propagate the column-A style (s=2) down through all the grand-total rows.
$q1.Range("A2").Copy()
$total.Range("A2:A4").PasteSpecial(-4122)

# Keep the first sheet as the active/selected tab, as in the original workbook.
$wb.Worksheets.Item(1).Activate()
